# This script applies a row-content permutation to the "Artfynd" sheet.
# Entire logical records (rows) were re-ordered/re-matched against their
# coordinates/species data; column A (Id) through AC (Publik kommentar)
# "follow" the record, while some rows gain/lose the optional
# K/L/M/N (Alder-Stadium/Kon/Aktivitet/Metod) and AC (Publik kommentar)
# cells depending on which record now occupies that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (the ORIGINAL content of the source row becomes
# the NEW content of the target row)
$mapping = @{
    2=3; 3=4; 4=2;
    5=6; 6=5;
    11=13; 13=11;
    14=15; 15=14;
    16=17; 17=16;
    18=19; 19=21; 21=18;
    23=24; 24=23;
    25=26; 26=27; 27=25;
    30=31; 31=30;
    32=35; 33=34; 34=33; 35=32;
}

# Columns whose value ever differs between records (everything else is
# identical boilerplate shared by every row, so it is left untouched).
$numericCols = @(1, 5, 17, 18)          # A, E, Q, R
$textCols    = @(2, 6, 7, 8, 26, 28)    # B, F, G, H, Z, AB

# rows that need the optional K/L/M/N/AC activity-block cells
$rowsNeedingActivity = @(19, 32)
# rows that need to lose the optional K/L/M/N/AC activity-block cells
$rowsLosingActivity = @(21, 35)

# ---------------------------------------------------------------------
# Step 1: snapshot every value we will need from the CURRENT (pre-edit)
# state, before any writes happen (several rows are both a source and a
# target, so everything must be captured up front).
# ---------------------------------------------------------------------
$snapshot = @{}
$rowsInvolved = @()
foreach ($k in $mapping.Keys) { $rowsInvolved += $k }
foreach ($v in $mapping.Values) { $rowsInvolved += $v }
$rowsInvolved = $rowsInvolved | Sort-Object -Unique

foreach ($r in $rowsInvolved) {
    $rowVals = @{}
    foreach ($c in $numericCols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    foreach ($c in $textCols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    # M (13) and AC (29) - activity / public comment text
    $rowVals[13] = $ws.Cells.Item($r, 13).Value2
    $rowVals[29] = $ws.Cells.Item($r, 29).Value2
    $snapshot[$r] = $rowVals
}

# ---------------------------------------------------------------------
# Step 2: write the new values into every target row.
# ---------------------------------------------------------------------
foreach ($target in ($mapping.Keys | Sort-Object)) {
    $source = $mapping[$target]
    $src = $snapshot[$source]

    foreach ($c in $numericCols) {
        $ws.Cells.Item($target, $c).Value = $src[$c]
    }
    foreach ($c in $textCols) {
        $ws.Cells.Item($target, $c).Value = $src[$c]
    }

    if ($rowsNeedingActivity -contains $target) {
        # this row now holds a "Tretåig hackspett" style record: it needs
        # K/L/N (blank) and M/AC (text) cells that did not exist before.
        $ws.Cells.Item($target, 11).Value = ""
        $ws.Cells.Item($target, 12).Value = ""
        $ws.Cells.Item($target, 13).Value = $src[13]
        $ws.Cells.Item($target, 14).Value = ""
        $ws.Cells.Item($target, 29).Value = $src[29]
    }
    elseif ($rowsLosingActivity -contains $target) {
        # this row now holds a record without an activity block: remove
        # the K/L/M/N/AC cells entirely.
        $ws.Cells.Item($target, 11).Value = $null
        $ws.Cells.Item($target, 12).Value = $null
        $ws.Cells.Item($target, 13).Value = $null
        $ws.Cells.Item($target, 14).Value = $null
        $ws.Cells.Item($target, 29).Value = $null
    }
    else {
        # M / AC keep existing presence, just update their text
        $ws.Cells.Item($target, 13).Value = $src[13]
        $ws.Cells.Item($target, 29).Value = $src[29]
    }
}

Write-Host "Row permutation applied."
